$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab4")

# ---------------------------------------------------------------------------
# Table 1 (ARRAYLIST) data block, rows 2-11 (B,C,D columns). Formulas are
# replaced with static values; cells beyond the new data set are cleared.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 640.625
$ws.Range("C2").Value = 671.875
$ws.Range("D2").Value = 31.25

$ws.Range("B3").Value = 2640.625
$ws.Range("C3").Value = 2625
$ws.Range("D3").Value = 78.125

$ws.Range("B4").Value = 10531.25
$ws.Range("C4").Value = 11906.25
$ws.Range("D4").Value = 187.5

$ws.Range("B5").Value = 45062.5
$ws.Range("C5").Value = 43140.625
$ws.Range("D5").Value = 437.5

$ws.Range("B6").Value = 183062.5
$ws.Range("C6").Value = 195375
$ws.Range("D6").Value = 1093.75

$ws.Range("B7").Value = 709890.625
$ws.Range("C7").Value = 769625
$ws.Range("D7").Value = 2609.375

$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = 5593.75

$ws.Range("B9").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = 14109.375

$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = 33000

$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()

# ---------------------------------------------------------------------------
# Table 2 (LINKED_LIST) data block, rows 15-24 (B,C,D columns).
# ---------------------------------------------------------------------------
$ws.Range("B15").Value = 55781.25
$ws.Range("C15").Value = 51640.625
$ws.Range("D15").Value = 2765.625

$ws.Range("B16").Value = 464328.125
$ws.Range("C16").Value = 417343.75
$ws.Range("D16").Value = 14046.875

$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("D17").Value = 65765.625

$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("D18").Value = 294265.625

$ws.Range("B19").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("D19").ClearContents()

$ws.Range("B20").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("D20").ClearContents()

$ws.Range("B21").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()

$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("D22").ClearContents()

$ws.Range("B23").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("D23").ClearContents()

$ws.Range("B24").ClearContents()
$ws.Range("C24").ClearContents()
$ws.Range("D24").ClearContents()

# ---------------------------------------------------------------------------
# Resize the two tables to match the shrunk data ranges.
# ---------------------------------------------------------------------------
$t1 = $ws.ListObjects.Item("Table1")
$t1.Resize($ws.Range("A1:D8"))

$t2 = $ws.ListObjects.Item("Table13")
$t2.Resize($ws.Range("A14:D16"))

# ---------------------------------------------------------------------------
# Selection / active sheet state to match the author's saved view.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("G9").Select()

Write-Host "Edits applied successfully"
